$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the "Não encontrado" registration placeholders with the brands'
# real registration numbers found by the new lookup function.
$ws.Range("D2").Value = 105830220
$ws.Range("D13").Value = 102350779
$ws.Range("D18").Value = 167730219
$ws.Range("D30").Value = 112360011
$ws.Range("D31").Value = 112360011
$ws.Range("D43").Value = 118190327
